# RunMode Config Update and Test Skip Functionality Added
#
# 1. Add a new "TestSuite" sheet at the front of the workbook containing a
#    TestCaseId / RunMode table used to decide which tests should run.
# 2. Rename the existing flow sheets to match the new test-suite naming
#    convention.
# 3. Style the header row (yellow fill + border) and data rows (border) on
#    all three sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet as the first tab -----------------------------
# NOTE: worksheet object references returned from $wb.Worksheets are
# positional, so any sheet-repositioning operation (like Add here) must
# happen before we grab the references we intend to keep using.
$testSuite = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$testSuite.Name = "TestSuite"

# --- 2. Rename the pre-existing sheets -------------------------------------
$addCustomer = $wb.Worksheets.Item("addCustomerFlow")
$addCustomer.Name = "AddCustomerTest"

$openAccount = $wb.Worksheets.Item("openAccountFlow")
$openAccount.Name = "OpenAccountTest"

# --- 3. Populate the TestSuite sheet ---------------------------------------
# Values are written in this particular column-major order so that the
# shared-string table ends up built in the same order as the source
# workbook (TestCaseId, RunMode, Y, N, BankManagerLoginTest,
# AddCustomerTest, OpenAccountTest).
$testSuite.Range("A1").Value = "TestCaseId"
$testSuite.Range("B1").Value = "RunMode"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("B4").Value = "N"
$testSuite.Range("A2").Value = "BankManagerLoginTest"
$testSuite.Range("A3").Value = "AddCustomerTest"
$testSuite.Range("A4").Value = "OpenAccountTest"

$testSuite.Columns.Item(1).ColumnWidth = 25.5

# --- 4. Apply styling -------------------------------------------------------
# Border-only style is created first (ends up as cellXfs index 1), then the
# yellow-fill + border header style (cellXfs index 2), matching the target
# workbook's style table ordering. The same two styles get reused (not
# recreated) on the other two sheets.
$testSuite.Range("A2:B4").Borders.LineStyle = 1
$testSuite.Range("A1:B1").Borders.LineStyle = 1
$testSuite.Range("A1:B1").Interior.Color = 65535

$addCustomer.Range("A2:D4").Borders.LineStyle = 1
$addCustomer.Range("A1:D1").Borders.LineStyle = 1
$addCustomer.Range("A1:D1").Interior.Color = 65535

$openAccount.Range("A2:C3").Borders.LineStyle = 1
$openAccount.Range("A1:C1").Borders.LineStyle = 1
$openAccount.Range("A1:C1").Interior.Color = 65535

# --- 5. Selections / active sheet ------------------------------------------
# Touch the other sheets' selection first, then finish on TestSuite so it
# ends up as the active/selected tab.
[void]$addCustomer.Activate()
[void]$addCustomer.Range("D17").Select()

[void]$openAccount.Activate()
[void]$openAccount.Range("D15").Select()

[void]$testSuite.Activate()
[void]$testSuite.Range("B4").Select()
